$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row: correct-answer marking value 3 -> 5
$ws.Range("B11").Value = 5

# Update "Total" row: total marks 75 -> 125
$ws.Range("B12").Value = 125

# Update the "Max" column total display "75/84" -> "125/140"
$ws.Range("E12").Value = "125/140"
